$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated imputed values for the KNN algorithm result sheet
$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.875
$ws.Range("D3").Value = -7.506
$ws.Range("D12").Value = -7.277000000000001
$ws.Range("A14").Value = -21.913
$ws.Range("A16").Value = -22.013
$ws.Range("C18").Value = -12.282
$ws.Range("A21").Value = -20.217
$ws.Range("A23").Value = -20.299
$ws.Range("C24").Value = -12.325
$ws.Range("D24").Value = -7.786
$ws.Range("A25").Value = -21.937
$ws.Range("C25").Value = -13.051
$ws.Range("D25").Value = -8.715
$ws.Range("A26").Value = -21.277
$ws.Range("C27").Value = -13.055
$ws.Range("A29").Value = -21.219
$ws.Range("C30").Value = -12.7
$ws.Range("C31").Value = -13.257
$ws.Range("C39").Value = -12.717
$ws.Range("A40").Value = -19.938
$ws.Range("D41").Value = -8.195000000000002
$ws.Range("C42").Value = -12.852
$ws.Range("C48").Value = -11.575
$ws.Range("D50").Value = -8.045999999999999
$ws.Range("C51").Value = -11.564
$ws.Range("C52").Value = -11.439
$ws.Range("A53").Value = -21.826
$ws.Range("D53").Value = -7.342999999999999
$ws.Range("C55").Value = -13.952
$ws.Range("C56").Value = -12.682
$ws.Range("D56").Value = -8.118
$ws.Range("A57").Value = -22.165
$ws.Range("C57").Value = -12.986
$ws.Range("D57").Value = -8.852
$ws.Range("D58").Value = -8.318000000000001
$ws.Range("A59").Value = -22.286
$ws.Range("C60").Value = -12.296
$ws.Range("D61").Value = -7.761999999999999
$ws.Range("D63").Value = -7.886
$ws.Range("D64").Value = -7.608
$ws.Range("A65").Value = -21.421
$ws.Range("A69").Value = -21.833
$ws.Range("D70").Value = -7.211
$ws.Range("D72").Value = -7.356999999999999
$ws.Range("C73").Value = -12.995
$ws.Range("C74").Value = -12.462
$ws.Range("A79").Value = -20.849
$ws.Range("A83").Value = -21.938
$ws.Range("D86").Value = -8.190999999999999
$ws.Range("C89").Value = -10.985
$ws.Range("D89").Value = -6.728
$ws.Range("C90").Value = -12.957
$ws.Range("A91").Value = -21.508
$ws.Range("C92").Value = -11.096
$ws.Range("A93").Value = -21.278
$ws.Range("D98").Value = -8.612
$ws.Range("A100").Value = -21.898
$ws.Range("D100").Value = -8.690999999999999
$ws.Range("D102").Value = -7.794000000000001
